$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Insert a new "Meta description" paragraph right after the document's
#    title (Heading1) paragraph:
#       Meta description: Enjoy Age of Halvar slot game with Wild Halvar and
#       bonus mode. Discover the pros and cons in our review and play for
#       free.
#    "Meta description" is bold, the rest of the sentence is plain text, and
#    the paragraph itself uses the default (Normal) style, not Heading1.
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Style = "Normal"

$metaRange = $metaPara.Range
$metaRange.Text = ": Enjoy Age of Halvar slot game with Wild Halvar and bonus mode. Discover the pros and cons in our review and play for free."

$labelInsertionPoint = $d.Range($metaRange.Start, $metaRange.Start)
$labelInsertionPoint.InsertBefore("Meta description")

$labelRange = $d.Range($metaRange.Start, $metaRange.Start + 16)
$labelRange.Bold = 1

# ---------------------------------------------------------------------------
# 2. Remove the duplicated bold title paragraph ("Play Age of Halvar Free:
#    Review and Features") that used to sit right before the closing
#    description paragraph near the end of the document.
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
$duplicateTitlePara = $d.Paragraphs.Item($count - 1)
$duplicateTitlePara.Range.Delete()

# ---------------------------------------------------------------------------
# 3. Replace the text of the final (italic) paragraph with the new prompt
#    text, keeping its existing italic formatting intact.
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
$closingPara = $d.Paragraphs.Item($count)
$closingRange = $closingPara.Range
$closingTextRange = $d.Range($closingRange.Start, $closingRange.End - 1)
$closingTextRange.Text = 'Create a feature image fitting the game "Age Of Halvar": - The image should be in cartoon style - The should feature a happy Maya warrior with glasses Sorry, there seems to be a confusion in your prompt. The game is called "Age Of Halvar", which is based on Vikings theme, but your prompt is asking for a feature image of a happy Maya warrior with glasses. Please clarify the prompt so I can provide an appropriate response.'
